$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 12: average of the "k" column (J) across the data rows.
$ws.Range("J12").Formula = "=AVERAGE(J2:J11)"

# Rows 14-17: summary statistics block - labels in column A, formulas in column B.
$ws.Range("A14").Value = "Average of SW(S*)/SW(OPT)"
$ws.Range("B14").Formula = "=AVERAGE(N2:N11)"

$ws.Range("A15").Value = "Average of SC(S*)/SC(OPT)"
$ws.Range("B15").Formula = "=AVERAGE(Z2:Z11)"

$ws.Range("A16").Value = "Worst of SW(S*)/SW(OPT)"
$ws.Range("B16").Formula = "=MIN(N2:N11)"

$ws.Range("A17").Value = "Worst of SC(S*)/SC(OPT)"
$ws.Range("B17").Formula = "=MAX(Z2:Z11)"

# Build the bold/size-12/vertically-centered style once on a scratch cell,
# then format-paint it onto B14:B17 in a single operation so only one new
# cell style record gets added (rather than one per property assignment).
$tmpl = $ws.Range("ZZ1")
$tmpl.Font.Bold = $true
$tmpl.Font.Size = 12
$tmpl.VerticalAlignment = -4108  # xlVAlignCenter

$tmpl.Copy()
$target = $ws.Range("B14:B17")
$target.PasteSpecial(-4122)  # xlPasteFormats
$tmpl.Clear()

# Row height for the new summary rows.
$ws.Range("A14:B17").RowHeight = 15.6

# Selection matching the post-edit workbook state.
$ws.Range("A14:B17").Select()

# Page setup gained an explicit paper size / orientation.
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1
